# Apply the benchmark-stats fix described in the commit message:
# "Fixed README.md stats and docx preparation for all DaCapo - JDK 21 - Z GC tests"
#
# The document is a single-column table where each row holds one stat value.
# A handful of rows get their value text replaced outright, and the last
# three rows (which originally packed several tab-separated numbers into one
# run) get collapsed down to a single short value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple value replacements (row index -> new text)
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "37"
    6  = "0.00006"
    8  = "0.00001"
    10 = "0.00004"
    12 = "0.00131"
    44 = "100"
    45 = "0"
    46 = "406"
}

foreach ($rowIndex in $updates.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $updates[$rowIndex]
}

Write-Output "applied"
